$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Crypto price/volume refresh (GitHub Actions data pull).
# Cells D/E (and the B/C coin name+link for two re-ranked rows) are plain
# text in the sheet (stored as inline strings), not numbers. Excel's COM
# layer auto-coerces a plain decimal-looking string (e.g. "0.9959") into a
# Number when assigned via .Value, so for those specific values we first pin
# the cell to the Text number format ("@") to preserve the original text type
# and formatting (e.g. trailing zeros). Values Excel would not mis-parse as a
# number (URLs, names, multi-dot price strings, padded percentages) are left
# on the default format, matching the source file exactly.

$ws.Range("D2").Value = '20.176.74'
$ws.Range("E2").Value = '  -1.05%  '

$ws.Range("D3").Value = '1.431.34'
$ws.Range("E3").Value = '  -0.49%  '

$ws.Range("E4").Value = '  -0.17%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.9959'
$ws.Range("E5").Value = '  -0.62%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '277.44'
$ws.Range("E6").Value = '  +0.03%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3710'
$ws.Range("E7").Value = '  -0.44%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3158'
$ws.Range("E8").Value = '  +2.33%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '40.28'
$ws.Range("E9").Value = '  -0.42%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.059'
$ws.Range("E10").Value = '  +4.69%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06587'

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9973'
$ws.Range("E12").Value = '  -0.51%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.562'
$ws.Range("E13").Value = '  +3.56%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '18.25'
$ws.Range("E14").Value = '  +5.38%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.238'
$ws.Range("E15").Value = '  +1.66%  '

$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001032'
$ws.Range("E16").Value = '  +1.95%  '

$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '1.433.66'
$ws.Range("E17").Value = '  -0.36%  '

$ws.Range("E18").Value = '  -1.28%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9959'
$ws.Range("E19").Value = '  -0.58%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '72.06'
$ws.Range("E20").Value = '  -5.89%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.631'
$ws.Range("E21").Value = '  -1.76%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '14.90'
$ws.Range("E22").Value = '  +3.39%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.17'
$ws.Range("E23").Value = '  +2.51%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.236'
$ws.Range("E24").Value = '  -3.61%  '

$ws.Range("D25").Value = '20.194.84'
$ws.Range("E25").Value = '  -0.96%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.310'
$ws.Range("E26").Value = '  +3.55%  '

$ws.Range("E27").Value = '  -5.46%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.48'
$ws.Range("E28").Value = '  +2.64%  '

$ws.Range("D29").Value = '1.590.90'
$ws.Range("E29").Value = '  -0.67%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '111.85'
$ws.Range("E30").Value = '  +1.51%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.952'
$ws.Range("E31").Value = '  +0.57%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.327'
$ws.Range("E32").Value = '  -2.84%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.8451'
$ws.Range("E33").Value = '  -7.53%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.07812'
$ws.Range("E34").Value = '  +1.38%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.490'
$ws.Range("E35").Value = '  +10.94%  '

$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.05912'
$ws.Range("E36").Value = '  +3.28%  '

$ws.Range("B37").Value = 'InternetComputer(DFINITY)'
$ws.Range("C37").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.936'
$ws.Range("E37").Value = '  +4.30%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '10.82'
$ws.Range("E38").Value = '  -1.31%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.9957'
$ws.Range("E39").Value = '  -0.60%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.856'
$ws.Range("E40").Value = '  -6.07%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.02071'
$ws.Range("E41").Value = '  +2.09%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.113'
$ws.Range("E42").Value = '  -2.06%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1886'
$ws.Range("E43").Value = '  -1.63%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5382'
$ws.Range("E44").Value = '  +1.03%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '12.39'
$ws.Range("E45").Value = '  +2.42%  '

$ws.Range("E46").Value = '  -0.72%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '119.35'
$ws.Range("E47").Value = '  +6.29%  '

$ws.Range("E48").Value = '  +2.30%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.802'
$ws.Range("E49").Value = '  +0.73%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.048'

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06274'
$ws.Range("E51").Value = '  -0.07%  '
